$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3
$ws.Cells.Item(3, 4).Value = 44453
$ws.Cells.Item(3, 8).Value = 'Camote'
$ws.Cells.Item(3, 9).Value = '1a nueva(o)'
$ws.Cells.Item(3, 10).Value = 800
$ws.Cells.Item(3, 11).Value = 630
$ws.Cells.Item(3, 12).Value = 650
$ws.Cells.Item(3, 13).Value = 640
$ws.Cells.Item(3, 15).Value = 'Perú'
$ws.Cells.Item(3, 16).Value = 640

# Row 4
$ws.Cells.Item(4, 4).Value = 44204
$ws.Cells.Item(4, 8).Value = 'Camote'
$ws.Cells.Item(4, 9).Value = '2a nueva(o)'
$ws.Cells.Item(4, 10).Value = 1600
$ws.Cells.Item(4, 11).Value = 500
$ws.Cells.Item(4, 12).Value = 550
$ws.Cells.Item(4, 13).Value = 525
$ws.Cells.Item(4, 15).Value = 'Región del Maule'
$ws.Cells.Item(4, 16).Value = 525

# Row 5
$ws.Cells.Item(5, 4).Value = 44211
$ws.Cells.Item(5, 8).Value = 'Camote'
$ws.Cells.Item(5, 9).Value = '1a nueva(o)'
$ws.Cells.Item(5, 10).Value = 1600
$ws.Cells.Item(5, 11).Value = 500
$ws.Cells.Item(5, 12).Value = 550
$ws.Cells.Item(5, 13).Value = 525
$ws.Cells.Item(5, 15).Value = 'Región de O''Higgins'
$ws.Cells.Item(5, 16).Value = 525

# Row 6
$ws.Cells.Item(6, 4).Value = 44399
$ws.Cells.Item(6, 8).Value = 'Camote'
$ws.Cells.Item(6, 9).Value = '1a (guarda)'
$ws.Cells.Item(6, 10).Value = 800
$ws.Cells.Item(6, 11).Value = 450
$ws.Cells.Item(6, 12).Value = 480
$ws.Cells.Item(6, 13).Value = 465
$ws.Cells.Item(6, 15).Value = 'Provincia de Melipilla'
$ws.Cells.Item(6, 16).Value = 465

# Row 7
$ws.Cells.Item(7, 4).Value = 44201
$ws.Cells.Item(7, 8).Value = 'Camote'
$ws.Cells.Item(7, 9).Value = '1a nueva(o)'
$ws.Cells.Item(7, 10).Value = 1360
$ws.Cells.Item(7, 11).Value = 730
$ws.Cells.Item(7, 12).Value = 750
$ws.Cells.Item(7, 13).Value = 740
$ws.Cells.Item(7, 15).Value = 'Perú'
$ws.Cells.Item(7, 16).Value = 740

# Row 8 is unchanged by this edit (left as-is)

# Row 9
$ws.Cells.Item(9, 4).Value = 44179
$ws.Cells.Item(9, 8).Value = 'Camote'
$ws.Cells.Item(9, 9).Value = '1a nueva(o)'
$ws.Cells.Item(9, 10).Value = 1000
$ws.Cells.Item(9, 11).Value = 1350
$ws.Cells.Item(9, 12).Value = 1400
$ws.Cells.Item(9, 13).Value = 1375
$ws.Cells.Item(9, 15).Value = 'Perú'
$ws.Cells.Item(9, 16).Value = 1375

# Row 10
$ws.Cells.Item(10, 4).Value = 44483
$ws.Cells.Item(10, 8).Value = 'Camote'
$ws.Cells.Item(10, 9).Value = '1a nueva(o)'
$ws.Cells.Item(10, 10).Value = 1300
$ws.Cells.Item(10, 11).Value = 550
$ws.Cells.Item(10, 12).Value = 580
$ws.Cells.Item(10, 13).Value = 565
$ws.Cells.Item(10, 15).Value = 'Perú'
$ws.Cells.Item(10, 16).Value = 565

# Row 11
$ws.Cells.Item(11, 4).Value = 44349
$ws.Cells.Item(11, 8).Value = 'Pachia'
$ws.Cells.Item(11, 9).Value = '1a nueva(o)'
$ws.Cells.Item(11, 10).Value = 1200
$ws.Cells.Item(11, 11).Value = 730
$ws.Cells.Item(11, 12).Value = 750
$ws.Cells.Item(11, 13).Value = 740
$ws.Cells.Item(11, 15).Value = 'Perú'
$ws.Cells.Item(11, 16).Value = 740

# Row 12
$ws.Cells.Item(12, 4).Value = 44243
$ws.Cells.Item(12, 8).Value = 'Camote'
$ws.Cells.Item(12, 9).Value = '2a nueva(o)'
$ws.Cells.Item(12, 10).Value = 1600
$ws.Cells.Item(12, 11).Value = 450
$ws.Cells.Item(12, 12).Value = 480
$ws.Cells.Item(12, 13).Value = 465
$ws.Cells.Item(12, 15).Value = 'Región del Maule'
$ws.Cells.Item(12, 16).Value = 465

# Row 13
$ws.Cells.Item(13, 4).Value = 44469
$ws.Cells.Item(13, 8).Value = 'Camote'
$ws.Cells.Item(13, 9).Value = '1a nueva(o)'
$ws.Cells.Item(13, 10).Value = 1200
$ws.Cells.Item(13, 11).Value = 600
$ws.Cells.Item(13, 12).Value = 650
$ws.Cells.Item(13, 13).Value = 625
$ws.Cells.Item(13, 15).Value = 'Perú'
$ws.Cells.Item(13, 16).Value = 625

# Row 14
$ws.Cells.Item(14, 4).Value = 44476
$ws.Cells.Item(14, 8).Value = 'Camote'
$ws.Cells.Item(14, 9).Value = '1a nueva(o)'
$ws.Cells.Item(14, 10).Value = 1200
$ws.Cells.Item(14, 11).Value = 480
$ws.Cells.Item(14, 12).Value = 500
$ws.Cells.Item(14, 13).Value = 490
$ws.Cells.Item(14, 15).Value = 'Perú'
$ws.Cells.Item(14, 16).Value = 490

# Row 15
$ws.Cells.Item(15, 4).Value = 44168
$ws.Cells.Item(15, 8).Value = 'Camote'
$ws.Cells.Item(15, 9).Value = '1a nueva(o)'
$ws.Cells.Item(15, 10).Value = 1200
$ws.Cells.Item(15, 11).Value = 1500
$ws.Cells.Item(15, 12).Value = 1700
$ws.Cells.Item(15, 13).Value = 1600
$ws.Cells.Item(15, 15).Value = 'Perú'
$ws.Cells.Item(15, 16).Value = 1600

# Row 16
$ws.Cells.Item(16, 4).Value = 44231
$ws.Cells.Item(16, 8).Value = 'Camote'
$ws.Cells.Item(16, 9).Value = '1a nueva(o)'
$ws.Cells.Item(16, 10).Value = 1300
$ws.Cells.Item(16, 11).Value = 450
$ws.Cells.Item(16, 12).Value = 480
$ws.Cells.Item(16, 13).Value = 465
$ws.Cells.Item(16, 15).Value = 'Perú'
$ws.Cells.Item(16, 16).Value = 465

# Row 17
$ws.Cells.Item(17, 4).Value = 44490
$ws.Cells.Item(17, 8).Value = 'Camote'
$ws.Cells.Item(17, 9).Value = '1a nueva(o)'
$ws.Cells.Item(17, 10).Value = 1200
$ws.Cells.Item(17, 11).Value = 450
$ws.Cells.Item(17, 12).Value = 480
$ws.Cells.Item(17, 13).Value = 465
$ws.Cells.Item(17, 15).Value = 'Perú'
$ws.Cells.Item(17, 16).Value = 465

# Row 18
$ws.Cells.Item(18, 4).Value = 44322
$ws.Cells.Item(18, 8).Value = 'Camote'
$ws.Cells.Item(18, 9).Value = '1a (cosecha)'
$ws.Cells.Item(18, 10).Value = 1200
$ws.Cells.Item(18, 11).Value = 350
$ws.Cells.Item(18, 12).Value = 400
$ws.Cells.Item(18, 13).Value = 375
$ws.Cells.Item(18, 15).Value = 'Región del Maule'
$ws.Cells.Item(18, 16).Value = 375

# Row 19
$ws.Cells.Item(19, 4).Value = 44238
$ws.Cells.Item(19, 8).Value = 'Camote'
$ws.Cells.Item(19, 9).Value = '1a nueva(o)'
$ws.Cells.Item(19, 10).Value = 1250
$ws.Cells.Item(19, 11).Value = 430
$ws.Cells.Item(19, 12).Value = 450
$ws.Cells.Item(19, 13).Value = 440
$ws.Cells.Item(19, 15).Value = 'Perú'
$ws.Cells.Item(19, 16).Value = 440

# Row 20
$ws.Cells.Item(20, 4).Value = 44175
$ws.Cells.Item(20, 8).Value = 'Camote'
$ws.Cells.Item(20, 9).Value = '1a nueva(o)'
$ws.Cells.Item(20, 10).Value = 1200
$ws.Cells.Item(20, 11).Value = 1500
$ws.Cells.Item(20, 12).Value = 1550
$ws.Cells.Item(20, 13).Value = 1525
$ws.Cells.Item(20, 15).Value = 'Perú'
$ws.Cells.Item(20, 16).Value = 1525

# Row 21
$ws.Cells.Item(21, 4).Value = 44428
$ws.Cells.Item(21, 8).Value = 'Camote'
$ws.Cells.Item(21, 9).Value = '1a nueva(o)'
$ws.Cells.Item(21, 10).Value = 600
$ws.Cells.Item(21, 11).Value = 580
$ws.Cells.Item(21, 12).Value = 600
$ws.Cells.Item(21, 13).Value = 590
$ws.Cells.Item(21, 15).Value = 'Perú'
$ws.Cells.Item(21, 16).Value = 590
